$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host ($ws | Get-Member | Out-String)
